$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09)
foreach ($row in 2..8) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
